$d = $word.ActiveDocument

$replacements = @(
    @("781×8=", "256×6="),
    @("542×8=", "937×4="),
    @("695×2=", "666×8="),
    @("167×6=", "681×5="),
    @("438×7=", "532×9="),
    @("328×5=", "852×2="),
    @("204×7=", "521×2="),
    @("793×4=", "221×9="),
    @("267×8=", "469×6="),
    @("370×6=", "897×4="),
    @("641×2=", "276×4="),
    @("168×4=", "219×6="),
    @("497×9=", "498×6="),
    @("408×8=", "190×7="),
    @("377×8=", "888×8="),
    @("583×9=", "669×4="),
    @("834×7=", "641×9="),
    @("778×5=", "251×9="),
    @("624×5=", "709×5="),
    @("176×6=", "650×6="),
    @("188×6=", "644×6="),
    @("438×9=", "218×9="),
    @("317×7=", "432×5="),
    @("585×8=", "236×2="),
    @("643×7=", "992×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
